# Apply updated cryptos list values (prices, 1h volume %, and a couple of row swaps)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.036.64'
$ws.Range('E2').Value = '  -3.40%  '
$ws.Range('D3').Value = '3.514.03'
$ws.Range('E3').Value = '  -4.31%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = "'611.97"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -5.11%  '
$ws.Range('D6').Value = "'154.06"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.65%  '
$ws.Range('D7').Value = '3.509.88'
$ws.Range('E7').Value = '  -4.41%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('E9').Value = '  -2.36%  '
$ws.Range('E10').Value = '  -2.99%  '
$ws.Range('D11').Value = "'6.84"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.47%  '
$ws.Range('D12').Value = "'0.429"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -4.67%  '
$ws.Range('D13').Value = "'0.0000222"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -4.78%  '
$ws.Range('D14').Value = "'32.09"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.00%  '
$ws.Range('D15').Value = '4.110.39'
$ws.Range('E15').Value = '  -4.19%  '
$ws.Range('D16').Value = '3.501.19'
$ws.Range('E16').Value = '  -4.44%  '
$ws.Range('D17').Value = '67.111.62'
$ws.Range('E17').Value = '  -3.32%  '
$ws.Range('E18').Value = '  +0.78%  '
$ws.Range('D19').Value = "'6.38"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.17%  '
$ws.Range('D20').Value = "'15.47"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.61%  '
$ws.Range('D21').Value = "'453.23"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.24%  '
$ws.Range('D22').Value = "'9.35"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -5.47%  '
$ws.Range('D23').Value = "'0.642"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.80%  '
$ws.Range('D24').Value = "'78.83"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.81%  '
$ws.Range('E25').Value = '  -0.16%  '
$ws.Range('D26').Value = '3.658.42'
$ws.Range('E26').Value = '  -4.18%  '
$ws.Range('E27').Value = '  -3.07%  '
$ws.Range('E28').Value = '  -4.35%  '
$ws.Range('D29').Value = "'8.30"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -8.64%  '
$ws.Range('E30').Value = '  -2.21%  '
$ws.Range('E31').Value = '  -2.96%  '
$ws.Range('E32').Value = '  +1.36%  '
$ws.Range('D33').Value = "'25.94"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.65%  '
$ws.Range('E34').Value = '  -5.45%  '
$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').Value = "'0.158"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.42%  '
$ws.Range('B36').Value = 'NEARProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D36').Value = "'6.19"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.18%  '
$ws.Range('D37').Value = '3.513.36'
$ws.Range('E37').Value = '  -4.03%  '
$ws.Range('D38').Value = "'8.03"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.86%  '
$ws.Range('E39').Value = '  +0.00%  '
$ws.Range('B40').Value = 'Monero'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D40').Value = "'178.16"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.39%  '
$ws.Range('B41').Value = 'FirstDigitalUSD'
$ws.Range('C41').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D41').Value = "'0.999"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.15%  '
$ws.Range('E42').Value = '  -5.43%  '
$ws.Range('E43').Value = '  -3.23%  '
$ws.Range('E44').Value = '  -3.15%  '
$ws.Range('E45').Value = '  -3.98%  '
$ws.Range('D46').Value = "'28.95"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +6.56%  '
$ws.Range('D47').Value = "'45.57"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.33%  '
$ws.Range('D48').Value = "'2.66"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.75%  '
$ws.Range('D49').Value = "'7.66"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.48%  '
$ws.Range('E50').Value = '  -3.22%  '
$ws.Range('E51').Value = '  -4.15%  '
